$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '24.897.84'
Set-TextValue $ws.Range("E2") '  +2.17%  '
Set-TextValue $ws.Range("D3") '1.711.88'
Set-TextValue $ws.Range("E3") '  +2.08%  '
Set-TextValue $ws.Range("E4") '  +0.01%  '
Set-TextValue $ws.Range("D5") '311.85'
Set-TextValue $ws.Range("E5") '  +1.73%  '
Set-TextValue $ws.Range("D6") '0.9999'
Set-TextValue $ws.Range("E6") '  +0.24%  '
Set-TextValue $ws.Range("E7") '  +1.20%  '
Set-TextValue $ws.Range("D8") '49.76'
Set-TextValue $ws.Range("E8") '  +3.40%  '
Set-TextValue $ws.Range("D9") '0.3454'
Set-TextValue $ws.Range("E9") '  +0.54%  '
Set-TextValue $ws.Range("D10") '1.209'
Set-TextValue $ws.Range("E10") '  +2.59%  '
Set-TextValue $ws.Range("D11") '0.07567'
Set-TextValue $ws.Range("E11") '  +4.23%  '
Set-TextValue $ws.Range("D12") '1.003'
Set-TextValue $ws.Range("E12") '  +0.29%  '
Set-TextValue $ws.Range("D13") '21.16'
Set-TextValue $ws.Range("E13") '  +4.12%  '
Set-TextValue $ws.Range("D14") '6.332'
Set-TextValue $ws.Range("E14") '  +3.85%  '
Set-TextValue $ws.Range("D15") '7.069'
Set-TextValue $ws.Range("E15") '  +4.84%  '
Set-TextValue $ws.Range("D16") '1.713.03'
Set-TextValue $ws.Range("E16") '  +1.96%  '
Set-TextValue $ws.Range("E17") '  +2.67%  '
Set-TextValue $ws.Range("E18") '  +0.20%  '
Set-TextValue $ws.Range("D19") '1.0000'
Set-TextValue $ws.Range("E19") '  +0.18%  '
Set-TextValue $ws.Range("E20") '  +4.78%  '
Set-TextValue $ws.Range("D21") '17.39'
Set-TextValue $ws.Range("E21") '  +5.90%  '
Set-TextValue $ws.Range("D22") '6.403'
Set-TextValue $ws.Range("D23") '13.25'
Set-TextValue $ws.Range("E23") '  +11.07%  '
Set-TextValue $ws.Range("D24") '24.845.00'
Set-TextValue $ws.Range("E24") '  +2.18%  '
Set-TextValue $ws.Range("D25") '2.455'
Set-TextValue $ws.Range("E25") '  +1.01%  '
Set-TextValue $ws.Range("D26") '2.801'
Set-TextValue $ws.Range("E26") '  +5.40%  '
Set-TextValue $ws.Range("D27") '20.47'
Set-TextValue $ws.Range("E27") '  +4.69%  '
Set-TextValue $ws.Range("D28") '152.20'
Set-TextValue $ws.Range("E28") '  +0.00%  '
Set-TextValue $ws.Range("D29") '132.99'
Set-TextValue $ws.Range("E29") '  +4.64%  '
Set-TextValue $ws.Range("D30") '1.903.81'
Set-TextValue $ws.Range("E30") '  +2.04%  '
Set-TextValue $ws.Range("D31") '1.243'
Set-TextValue $ws.Range("E31") '  +28.77%  '
Set-TextValue $ws.Range("D32") '6.942'
Set-TextValue $ws.Range("E32") '  +9.89%  '
Set-TextValue $ws.Range("D33") '4.237'
Set-TextValue $ws.Range("E33") '  +5.26%  '
Set-TextValue $ws.Range("D34") '1.851'
Set-TextValue $ws.Range("E34") '  +6.23%  '
Set-TextValue $ws.Range("D35") '13.86'
Set-TextValue $ws.Range("E35") '  +12.81%  '
Set-TextValue $ws.Range("D36") '0.08807'
Set-TextValue $ws.Range("E36") '  +3.90%  '
Set-TextValue $ws.Range("D37") '5.629'
Set-TextValue $ws.Range("E37") '  +5.56%  '
Set-TextValue $ws.Range("D38") '0.06709'
Set-TextValue $ws.Range("E38") '  +3.87%  '
Set-TextValue $ws.Range("D39") '9.332'
Set-TextValue $ws.Range("E39") '  +3.26%  '
Set-TextValue $ws.Range("D40") '0.02419'
Set-TextValue $ws.Range("E40") '  +3.92%  '
Set-TextValue $ws.Range("D41") '0.2244'
Set-TextValue $ws.Range("E41") '  +6.53%  '
Set-TextValue $ws.Range("D42") '1.280'
Set-TextValue $ws.Range("E42") '  +1.56%  '
Set-TextValue $ws.Range("D43") '0.6476'
Set-TextValue $ws.Range("E43") '  +5.03%  '
Set-TextValue $ws.Range("D44") '0.9996'
Set-TextValue $ws.Range("E44") '  +0.20%  '
Set-TextValue $ws.Range("D45") '14.05'
Set-TextValue $ws.Range("E45") '  +7.63%  '
Set-TextValue $ws.Range("D46") '0.6194'
Set-TextValue $ws.Range("E46") '  +4.21%  '
Set-TextValue $ws.Range("E47") '  +1.29%  '
Set-TextValue $ws.Range("D48") '2.145'
Set-TextValue $ws.Range("E48") '  +6.06%  '
Set-TextValue $ws.Range("D49") '130.35'
Set-TextValue $ws.Range("E49") '  +2.55%  '
Set-TextValue $ws.Range("D50") '0.07333'
Set-TextValue $ws.Range("E50") '  +1.63%  '
Set-TextValue $ws.Range("D51") '80.17'
Set-TextValue $ws.Range("E51") '  +5.77%  '
